$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 137, shifting the existing data (rows 137-192)
# down to rows 138-193.
$ws.Range("A137:R137").Insert(-4121)

# Populate the newly inserted row 137 with a new weekly record (same
# reference data as the row that used to be at 137, now at 138, but with
# an updated date and volume).
$ws.Range("A137").Value = 10
$ws.Range("B137").Value = "Vega Modelo de Temuco"
$ws.Range("C137").Value = "La Araucanía"
$ws.Range("D137").Value = 44825
$ws.Range("E137").Value = 9
$ws.Range("F137").Value = 100114007
$ws.Range("G137").Value = "Jengibre"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 80
$ws.Range("K137").Value = 20000
$ws.Range("L137").Value = 20000
$ws.Range("M137").Value = 20000
$ws.Range("N137").Value = "$/caja 13 kilos"
$ws.Range("O137").Value = "Perú"
$ws.Range("P137").Value = 1538
$ws.Range("Q137").Value = 13
$ws.Range("R137").Value = "Hortaliza"
